$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Remove the "Santa Fe Springs, CA 90670" row (row 8) entirely and
#    shift the rows below it up. This also drops the now-unused
#    shared string and renumbers everything automatically.
# ------------------------------------------------------------------
$ws.Range("A8:P8").Delete(-4162) | Out-Null   # -4162 = xlShiftUp

# ------------------------------------------------------------------
# 2) Update the numeric data for rows 2-19 (after the shift above).
#    Helper to set a value with the same "0.00" number format used
#    by the rest of the data cells (style index 1 in the original
#    workbook).
# ------------------------------------------------------------------
function Set-DataCell($addr, $value) {
    $c = $ws.Range($addr)
    $c.Value = $value
    $c.NumberFormat = "0.00"
}

function Clear-DataCell($addr) {
    $ws.Range($addr).ClearContents() | Out-Null
}

# Row 2 - Olympia, WA 98516
Set-DataCell "B2" 440
Set-DataCell "C2" 830
Set-DataCell "D2" 1190

# Row 3 - Irving, TX 75061
Set-DataCell "B3" 425
Set-DataCell "C3" 805
Set-DataCell "D3" 1340

# Row 4 - Minooka, IL 60447
Set-DataCell "B4" 450
Set-DataCell "C4" 840
Set-DataCell "D4" 1190

# Row 5 - Nazareth, PA 18064
Set-DataCell "B5" 515
Set-DataCell "C5" 955
Set-DataCell "D5" 1315

# Row 6 - Bloomfield, CT 06002
Set-DataCell "B6" 570
Set-DataCell "C6" 1085
Set-DataCell "D6" 1535

# Row 7 - Daytona Beach, FL 32117
Set-DataCell "B7" 655
Set-DataCell "C7" 1245
Set-DataCell "D7" 1790

# Row 8 - Tracy, CA 95304 (was row 9 before the shift)
Set-DataCell "B8" 395
Set-DataCell "C8" 530
Set-DataCell "D8" 665
Set-DataCell "E8" 785
Set-DataCell "F8" 880
Clear-DataCell "G8"
Clear-DataCell "H8"
Clear-DataCell "I8"
Clear-DataCell "J8"
Clear-DataCell "K8"
Clear-DataCell "L8"
Clear-DataCell "M8"
Clear-DataCell "N8"
Clear-DataCell "O8"
Clear-DataCell "P8"

# Row 9 - Fort Worth, TX 76140 (was row 10 before the shift)
Set-DataCell "B9" 485
Clear-DataCell "C9"
Clear-DataCell "D9"
Clear-DataCell "E9"

# Row 10 - Modesto, CA 95353 (was row 11 before the shift)
Set-DataCell "B10" 395
Set-DataCell "C10" 530
Set-DataCell "D10" 665
Set-DataCell "E10" 785
Set-DataCell "F10" 880
Set-DataCell "G10" 980

# Row 11 - Dallas, TX 75244 (was row 12 before the shift)
Clear-DataCell "B11"
Clear-DataCell "C11"
Clear-DataCell "D11"
Clear-DataCell "E11"
Clear-DataCell "F11"
Set-DataCell "G11" 1800

# Row 12 - Chicago, IL 60628 (was row 13 before the shift)
Set-DataCell "E12" 1435
Clear-DataCell "G12"

# Row 13 - Stockton, CA 95205 (was row 14 before the shift)
Set-DataCell "B13" 335
Set-DataCell "C13" 480
Set-DataCell "D13" 610
Set-DataCell "E13" 700
Set-DataCell "F13" 840
Set-DataCell "G13" 975

# Row 14 - Stockton, CA 96215 (was row 15 before the shift)
Set-DataCell "B14" 400
Set-DataCell "C14" 540
Set-DataCell "D14" 690
Set-DataCell "E14" 820
Set-DataCell "F14" 910
Set-DataCell "G14" 1000
Set-DataCell "H14" 1050
Set-DataCell "I14" 1105
Set-DataCell "J14" 1150
Set-DataCell "K14" 1200

# Row 15 - Fremont, CA 94538 (was row 16 before the shift)
Set-DataCell "B15" 395
Set-DataCell "C15" 530
Set-DataCell "D15" 665
Set-DataCell "E15" 785
Set-DataCell "F15" 880
Set-DataCell "G15" 980
Clear-DataCell "H15"
Clear-DataCell "I15"
Clear-DataCell "J15"
Clear-DataCell "K15"

# Row 16 - Pleasant Grove, CA 95668 (was row 17 before the shift)
Set-DataCell "B16" 395
Set-DataCell "C16" 530
Set-DataCell "D16" 665
Set-DataCell "E16" 785
Set-DataCell "F16" 880
Set-DataCell "G16" 980

# Row 17 - Fairfield, CA 94533 (was row 18 before the shift)
Set-DataCell "B17" 400
Set-DataCell "C17" 565
Set-DataCell "D17" 705
Set-DataCell "E17" 835
Set-DataCell "F17" 930
Set-DataCell "G17" 1030

# Row 18 - Manteca, CA 95336 (was row 19 before the shift)
Set-DataCell "B18" 400
Set-DataCell "C18" 565
Set-DataCell "D18" 705
Set-DataCell "E18" 835
Set-DataCell "F18" 930
Set-DataCell "G18" 1030

# Row 19 - Milpitas, CA 95035 (was row 20 before the shift)
Set-DataCell "B19" 400
Set-DataCell "C19" 565
Set-DataCell "D19" 705
Set-DataCell "E19" 835
Set-DataCell "F19" 930
Set-DataCell "G19" 1030

# ------------------------------------------------------------------
# 3) Row 20 loses every value except a new one in column C, and a
#    brand-new row 21 is appended for "Loveland, CO 80538".
# ------------------------------------------------------------------
$ws.Range("A20:P20").Clear() | Out-Null
$ws.Range("A20").Value = "Boise, ID 83717"
Set-DataCell "C20" 1195

$ws.Range("A21").Value = "Loveland, CO 80538"
Set-DataCell "B21" 540

# ------------------------------------------------------------------
# 4) Cosmetic sheet-level changes: active cell/selection and the
#    width of column A.
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 23.5703125
$ws.Range("K7").Select() | Out-Null
